# Extend the yearly data table in column R (year 2021) to match column Q's
# formatting, then populate the header + data values for rows 3-33.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (number format, font, borders, etc.) from the existing
# "2020" column (Q) into the new "2021" column (R) for the whole table body.
$ws.Range("Q3:Q33").Copy()
$ws.Range("R3:R33").PasteSpecial(-4122)

# Header
$ws.Range("R3").Value = 2021

# Data rows
$ws.Range("R4").Value = 1.7931687443515183
$ws.Range("R5").Value = 1.0977143806517458
$ws.Range("R6").Value = 2.4989281705678046
$ws.Range("R7").Value = 2.3489023398681002
$ws.Range("R8").Value = 1.8410239038543676
$ws.Range("R9").Value = 2.8382683724659588
$ws.Range("R10").Value = 1.2584206034913306
$ws.Range("R11").Value = 0.79202525610136665
$ws.Range("R12").Value = 1.7183687369364922
$ws.Range("R13").Value = 1.7860084101151579
$ws.Range("R14").Value = 1.5807090270340762
$ws.Range("R15").Value = 1.9930959157478496
$ws.Range("R16").Value = 1.0231016349164126
$ws.Range("R17").Value = 0
$ws.Range("R18").Value = 2.0091214112068791
$ws.Range("R19").Value = 2.2092990108041848
$ws.Range("R20").Value = 0.86496336159360854
$ws.Range("R21").Value = 3.5236628052020538
$ws.Range("R22").Value = 1.4678252700798498
$ws.Range("R23").Value = 0.74155920237892192
$ws.Range("R24").Value = 2.1792664589099311
$ws.Range("R25").Value = 1.5302890103825006
$ws.Range("R26").Value = 0.80351618683358383
$ws.Range("R27").Value = 2.280288974802807
$ws.Range("R28").Value = 2.3014726663297309
$ws.Range("R29").Value = 1.7358308467556451
$ws.Range("R30").Value = 2.9402079315049163
$ws.Range("R31").Value = 1.2198989923634325
$ws.Range("R32").Value = 1.1878318505232399
$ws.Range("R33").Value = 1.2537455648750642

# Match the saved selection state recorded in the workbook.
[void]$ws.Range("S14").Select()
